$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the shared-string table gets recompacted
# (drops the now-unused RFB/RG1/RG2/eev strings) the same way the
# authoritative workbook does.
$ws.Cells.Clear()

# --- Header row ---
$ws.Range("A1").Value = "Inputs "
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Value = "This is the TI app note: SLAA869"
$ws.Range("F1").Value = "Inputs "
$ws.Range("F1").Font.Bold = $true
$ws.Range("G1").Value = "This is my design"

# --- Row 5 : VREF ---
$ws.Range("A5").Value = "VREF"
$ws.Range("B5").Value = 2.5
$ws.Range("B5").NumberFormat = "0.00E+00"
$ws.Range("F5").Value = "VREF"
$ws.Range("G5").Value = 2.5
$ws.Range("G5").NumberFormat = "0.00E+00"

# --- Row 2 col C/H : Vout1 ---
$ws.Range("C2").Value = "Vout1"
$ws.Range("C2").Font.Bold = $true
$ws.Range("H2").Value = "Vout1"
$ws.Range("H2").Font.Bold = $true

# --- Row 3 col C/H : Vout2 ---
$ws.Range("C3").Value = "Vout2"
$ws.Range("C3").Font.Bold = $true
$ws.Range("H3").Value = "Vout2"
$ws.Range("H3").Font.Bold = $true

# --- Row 4 col C/H : Vout3 ---
$ws.Range("C4").Value = "Vout3"
$ws.Range("C4").Font.Bold = $true
$ws.Range("H4").Value = "Vout3"
$ws.Range("H4").Font.Bold = $true

# --- Row 6-8 : VDAC ---
$ws.Range("A6").Value = "VDAC"
$ws.Range("B6").Value = 0
$ws.Range("B6").NumberFormat = "0.00E+00"
$ws.Range("F6").Value = "VDAC"
$ws.Range("G6").Value = 0
$ws.Range("G6").NumberFormat = "0.00E+00"

$ws.Range("A7").Value = "VDAC"
$ws.Range("B7").Value = 2.5
$ws.Range("B7").NumberFormat = "0.00E+00"
$ws.Range("F7").Value = "VDAC"
$ws.Range("G7").Value = 2.5
$ws.Range("G7").NumberFormat = "0.00E+00"

$ws.Range("A8").Value = "VDAC"
$ws.Range("B8").Value = 1.25
$ws.Range("B8").NumberFormat = "0.00E+00"
$ws.Range("F8").Value = "VDAC"
$ws.Range("G8").Value = 1.25
$ws.Range("G8").NumberFormat = "0.00E+00"

# --- Note below the table (first brand-new string) ---
$ws.Range("A21").Value = "see Figure 78 and equation 4 of DACx0508 datasheet"

# --- Row 2 : R3 / Vout1 ---
$ws.Range("A2").Value = "R3"
$ws.Range("B2").Value = 30000
$ws.Range("B2").NumberFormat = "0.00E+00"
$ws.Range("D2").Formula = "=(1 + B`$2/B`$3 + B`$2/B`$4)*B6-B`$2/B`$3*B`$5"
$ws.Range("D2").NumberFormat = "0.00E+00"

$ws.Range("F2").Value = "R3"
$ws.Range("G2").Value = 15000
$ws.Range("G2").NumberFormat = "0.00E+00"
$ws.Range("I2").Formula = "=(1 + G`$2/G`$3 + G`$2/G`$4)*G6-G`$2/G`$3*G`$5"
$ws.Range("I2").NumberFormat = "0.00E+00"

# --- Row 3 : R1 / Vout2 ---
$ws.Range("A3").Value = "R1"
$ws.Range("B3").Value = 7500
$ws.Range("B3").NumberFormat = "0.00E+00"
$ws.Range("D3").Formula = "=(1 + B`$2/B`$3 + B`$2/B`$4)*B7-B`$2/B`$3*B`$5"
$ws.Range("D3").NumberFormat = "0.00E+00"

$ws.Range("F3").Value = "R1"
$ws.Range("G3").Value = 10000
$ws.Range("G3").NumberFormat = "0.00E+00"
$ws.Range("I3").Formula = "=(1 + G`$2/G`$3 + G`$2/G`$4)*G7-G`$2/G`$3*G`$5"
$ws.Range("I3").NumberFormat = "0.00E+00"

# --- Row 4 : R2 / Vout3 ---
$ws.Range("A4").Value = "R2"
$ws.Range("B4").Value = 10000
$ws.Range("B4").NumberFormat = "0.00E+00"
$ws.Range("D4").Formula = "=(1 + B`$2/B`$3 + B`$2/B`$4)*B8-B`$2/B`$3*B`$5"
$ws.Range("D4").NumberFormat = "0.00E+00"

$ws.Range("F4").Value = "R2"
$ws.Range("G4").Value = 30000
$ws.Range("G4").NumberFormat = "0.00E+00"
$ws.Range("I4").Formula = "=(1 + G`$2/G`$3 + G`$2/G`$4)*G8-G`$2/G`$3*G`$5"
$ws.Range("I4").NumberFormat = "0.00E+00"

# --- Footprint note ---
$ws.Range("G23").Value = "RC0603FR-07240RL"

# --- Vstep per code (last brand-new string) ---
$ws.Range("K2").Value = "Vstep per code"
$ws.Range("L2").Formula = "=(I3-I2)/2^16"
$ws.Range("L2").NumberFormat = "0.00E+00"

# --- Column K width ---
$ws.Columns.Item(11).ColumnWidth = 13.5

# --- Selection / sheet name ---
$ws.Range("M2").Select() | Out-Null
$ws.Name = "dac80508_bipolar"
